$d = $word.ActiveDocument

$d.Content.Find.Execute("720×4=", $true, $false, $false, $false, $false, $true, 1, $false, "474×7=", 2) | Out-Null
$d.Content.Find.Execute("921×2=", $true, $false, $false, $false, $false, $true, 1, $false, "297×5=", 2) | Out-Null
$d.Content.Find.Execute("113×9=", $true, $false, $false, $false, $false, $true, 1, $false, "742×2=", 2) | Out-Null
$d.Content.Find.Execute("563×5=", $true, $false, $false, $false, $false, $true, 1, $false, "282×5=", 2) | Out-Null
$d.Content.Find.Execute("930×3=", $true, $false, $false, $false, $false, $true, 1, $false, "473×7=", 2) | Out-Null
$d.Content.Find.Execute("469×7=", $true, $false, $false, $false, $false, $true, 1, $false, "615×4=", 2) | Out-Null
$d.Content.Find.Execute("747×5=", $true, $false, $false, $false, $false, $true, 1, $false, "256×3=", 2) | Out-Null
$d.Content.Find.Execute("280×5=", $true, $false, $false, $false, $false, $true, 1, $false, "775×2=", 2) | Out-Null
$d.Content.Find.Execute("523×5=", $true, $false, $false, $false, $false, $true, 1, $false, "646×2=", 2) | Out-Null
$d.Content.Find.Execute("484×5=", $true, $false, $false, $false, $false, $true, 1, $false, "970×8=", 2) | Out-Null
$d.Content.Find.Execute("756×7=", $true, $false, $false, $false, $false, $true, 1, $false, "916×4=", 2) | Out-Null
$d.Content.Find.Execute("331×6=", $true, $false, $false, $false, $false, $true, 1, $false, "324×2=", 2) | Out-Null
$d.Content.Find.Execute("536×3=", $true, $false, $false, $false, $false, $true, 1, $false, "352×8=", 2) | Out-Null
$d.Content.Find.Execute("730×6=", $true, $false, $false, $false, $false, $true, 1, $false, "775×9=", 2) | Out-Null
$d.Content.Find.Execute("910×9=", $true, $false, $false, $false, $false, $true, 1, $false, "985×8=", 2) | Out-Null
$d.Content.Find.Execute("406×5=", $true, $false, $false, $false, $false, $true, 1, $false, "898×9=", 2) | Out-Null
$d.Content.Find.Execute("400×5=", $true, $false, $false, $false, $false, $true, 1, $false, "304×5=", 2) | Out-Null
$d.Content.Find.Execute("172×9=", $true, $false, $false, $false, $false, $true, 1, $false, "444×4=", 2) | Out-Null
$d.Content.Find.Execute("122×4=", $true, $false, $false, $false, $false, $true, 1, $false, "984×4=", 2) | Out-Null
$d.Content.Find.Execute("389×3=", $true, $false, $false, $false, $false, $true, 1, $false, "594×8=", 2) | Out-Null
$d.Content.Find.Execute("495×7=", $true, $false, $false, $false, $false, $true, 1, $false, "648×3=", 2) | Out-Null
$d.Content.Find.Execute("716×9=", $true, $false, $false, $false, $false, $true, 1, $false, "357×7=", 2) | Out-Null
$d.Content.Find.Execute("696×3=", $true, $false, $false, $false, $false, $true, 1, $false, "702×3=", 2) | Out-Null
$d.Content.Find.Execute("596×3=", $true, $false, $false, $false, $false, $true, 1, $false, "244×9=", 2) | Out-Null
$d.Content.Find.Execute("911×3=", $true, $false, $false, $false, $false, $true, 1, $false, "225×5=", 2) | Out-Null
